# DAI Municipios Enlace Transparencia.xlsx - apply the update described by the
# commit "Update DAI Municipios Enlace Transparencia.xlsx": fill in response
# details (Solicitud_Code, Fecha_Respuesta, Profesión Enlace Transparencia,
# Antigüedad en el cargo, Modalidad contratación, and an Obs. note) for a
# handful of "Codes" sheet rows that previously had empty C:H cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Helper: write a date value into a cell while re-using the date-formatted
# style (xf 12) that the sheet already applies to other D/E cells, instead of
# letting Excel synthesize a brand-new number-formatted style.
function Set-DateCell {
    param($row, $col, $serial)
    # D80 already carries the date style used throughout column D/E.
    $ws.Cells.Item(80, 4).Copy() | Out-Null
    $ws.Cells.Item($row, $col).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Cells.Item($row, $col).Value = $serial
}

$excel.CutCopyMode = $false

# ---- Row 44 (Comuna MU115) ----
$ws.Cells.Item(44, 3).Value = "CT001T0015682"
Set-DateCell 44 4 44659
Set-DateCell 44 5 44673
$ws.Cells.Item(44, 6).Value = "Secretaria Ejecutiva Bilingüe"
$ws.Cells.Item(44, 7).Value = "9 años y 10 meses"
$ws.Cells.Item(44, 8).Value = "Planta"

# ---- Row 80 (Comuna MU100) ----
Set-DateCell 80 5 44673
$ws.Cells.Item(80, 6).Value = "Ingenieria en Administración de Empresas"
$ws.Cells.Item(80, 7).Value = "6 años"
$ws.Cells.Item(80, 8).Value = "Contrata"
$ws.Rows.Item(80).RowHeight = 28.8

# ---- Row 107 (Comuna MU305) ----
$ws.Cells.Item(107, 3).Value = "CT001T0015682"
Set-DateCell 107 4 44643
Set-DateCell 107 5 44673
$ws.Cells.Item(107, 6).Value = "No informa (apoyo de informática)"
$ws.Cells.Item(107, 7).Value = "7 años"
$ws.Cells.Item(107, 8).Value = "Código del Trabajo"
$ws.Rows.Item(107).RowHeight = 28.8

# ---- Row 113 (Comuna MU206) ----
Set-DateCell 113 5 44673
$ws.Cells.Item(113, 6).Value = "Técnico Asistente Judicial"
$ws.Cells.Item(113, 7).Value = "2 años"
$ws.Cells.Item(113, 8).Value = "Contrata"

# ---- Row 171 (Comuna MU136) ----
Set-DateCell 171 5 44673
$ws.Cells.Item(171, 6).Value = "Ingeniero en Prevención de Riesgos (Transparencia Activa)"
$ws.Cells.Item(171, 7).Value = "3 años"
$ws.Cells.Item(171, 8).Value = "Contrata"
# I171 is a brand-new cell (row previously only spanned A:H); give it the
# same centered/wrapped style (xf 5) used by the other Obs.-style cells, e.g.
# I25, before writing its text.
$ws.Cells.Item(25, 9).Copy() | Out-Null
$ws.Cells.Item(171, 9).PasteSpecial($xlPasteFormats) | Out-Null
$ws.Cells.Item(171, 9).Value = "Además se menciona a un Técnico Nivel Medio Contador General (Planta) con 1 año y 8 meses de antigüedad en el cargo cumpliendo funciones relativas a Transparencia Pasiva "
$ws.Rows.Item(171).RowHeight = 57.6

# ---- Row 188 (Comuna MU088) ----
Set-DateCell 188 5 44673
$ws.Cells.Item(188, 6).Value = "Ingeniero en Administración de Empresas"
$ws.Cells.Item(188, 7).Value = "No informa"
$ws.Cells.Item(188, 8).Value = "Planta"
$ws.Rows.Item(188).RowHeight = 28.8

# ---- Row 216 (Comuna MU235) ----
$ws.Cells.Item(216, 3).Value = "CT001T0015682"
Set-DateCell 216 4 44659
Set-DateCell 216 5 44673
$ws.Cells.Item(216, 6).Value = "Ingeniero de Ejecución en Administración"
$ws.Cells.Item(216, 7).Value = "2 años"
$ws.Cells.Item(216, 8).Value = "Planta"
$ws.Rows.Item(216).RowHeight = 28.8

# ---- Row 333 (Comuna MU139) ----
$ws.Cells.Item(333, 3).Value = "CT001T0015682"
Set-DateCell 333 4 44643
Set-DateCell 333 5 44673
$ws.Cells.Item(333, 6).Value = "Abogado"
$ws.Cells.Item(333, 7).Value = "5 años"
$ws.Cells.Item(333, 8).Value = "Planta"

$excel.CutCopyMode = $false

# Final cursor position the author left the sheet in.
$ws.Range("C7").Select() | Out-Null
